$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $val) {
    $c = $ws.Range($addr)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = $origStyle
}

Set-TextCell $ws "B2" 'Bitcoin'
Set-TextCell $ws "C2" 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
Set-TextCell $ws "D2" '25.974.85'
Set-TextCell $ws "E2" '  +0.14%  '

Set-TextCell $ws "B3" 'Ethereum'
Set-TextCell $ws "C3" 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
Set-TextCell $ws "D3" '1.741.69'
Set-TextCell $ws "E3" '  +0.02%  '

Set-TextCell $ws "B4" 'TetherUSD'
Set-TextCell $ws "C4" 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
Set-TextCell $ws "D4" '0.9989'
Set-TextCell $ws "E4" '  -0.07%  '

Set-TextCell $ws "B5" 'BNB'
Set-TextCell $ws "C5" 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
Set-TextCell $ws "D5" '247.64'
Set-TextCell $ws "E5" '  +6.92%  '

Set-TextCell $ws "B6" 'USDC'
Set-TextCell $ws "C6" 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
Set-TextCell $ws "D6" '0.9989'
Set-TextCell $ws "E6" '  -0.11%  '

Set-TextCell $ws "B7" 'XRP'
Set-TextCell $ws "C7" 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
Set-TextCell $ws "D7" '0.5125'
Set-TextCell $ws "E7" '  -2.27%  '

Set-TextCell $ws "B8" 'Cardano'
Set-TextCell $ws "C8" 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextCell $ws "D8" '0.2755'
Set-TextCell $ws "E8" '  -0.18%  '

Set-TextCell $ws "B9" 'Dogecoin'
Set-TextCell $ws "C9" 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextCell $ws "D9" '0.06183'
Set-TextCell $ws "E9" '  +0.34%  '

Set-TextCell $ws "B10" 'WrappedEther'
Set-TextCell $ws "C10" 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextCell $ws "D10" '1.741.26'
Set-TextCell $ws "E10" '  -0.05%  '

Set-TextCell $ws "B11" 'TRON'
Set-TextCell $ws "C11" 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextCell $ws "D11" '0.07207'
Set-TextCell $ws "E11" '  +0.94%  '

Set-TextCell $ws "B12" 'Solana'
Set-TextCell $ws "C12" 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextCell $ws "D12" '15.00'
Set-TextCell $ws "E12" '  -1.49%  '

Set-TextCell $ws "B13" 'Polygon'
Set-TextCell $ws "C13" 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextCell $ws "D13" '0.6431'
Set-TextCell $ws "E13" '  -0.09%  '

Set-TextCell $ws "B14" 'Polkadot'
Set-TextCell $ws "C14" 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextCell $ws "D14" '4.628'
Set-TextCell $ws "E14" '  +2.05%  '

Set-TextCell $ws "B15" 'Litecoin'
Set-TextCell $ws "C15" 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextCell $ws "D15" '77.33'
Set-TextCell $ws "E15" '  -0.15%  '

Set-TextCell $ws "B16" 'Dai'
Set-TextCell $ws "C16" 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextCell $ws "D16" '0.9995'
Set-TextCell $ws "E16" '  +0.00%  '

Set-TextCell $ws "B17" 'BinanceUSD'
Set-TextCell $ws "C17" 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextCell $ws "D17" '0.9997'
Set-TextCell $ws "E17" '  +0.05%  '

Set-TextCell $ws "B18" 'WrappedBTC'
Set-TextCell $ws "C18" 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextCell $ws "D18" '25.985.01'
Set-TextCell $ws "E18" '  +0.31%  '

Set-TextCell $ws "B19" 'Avalanche'
Set-TextCell $ws "C19" 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextCell $ws "D19" '11.76'
Set-TextCell $ws "E19" '  +1.98%  '

Set-TextCell $ws "B20" 'ShibaInu'
Set-TextCell $ws "C20" 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextCell $ws "D20" '0.000006797'
Set-TextCell $ws "E20" '  +1.91%  '

Set-TextCell $ws "B21" 'WrappedliquidstakedEther2.0'
Set-TextCell $ws "C21" 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextCell $ws "D21" '1.962.36'
Set-TextCell $ws "E21" '  -0.05%  '

Set-TextCell $ws "B22" 'Uniswap'
Set-TextCell $ws "C22" 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextCell $ws "D22" '4.256'
Set-TextCell $ws "E22" '  -0.79%  '

Set-TextCell $ws "B23" 'Cosmos'
Set-TextCell $ws "C23" 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextCell $ws "D23" '8.622'
Set-TextCell $ws "E23" '  -1.46%  '

Set-TextCell $ws "B24" 'Chainlink'
Set-TextCell $ws "C24" 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextCell $ws "D24" '5.301'
Set-TextCell $ws "E24" '  +2.80%  '

Set-TextCell $ws "B25" 'Monero'
Set-TextCell $ws "C25" 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextCell $ws "D25" '136.60'
Set-TextCell $ws "E25" '  -2.94%  '

Set-TextCell $ws "B26" 'Toncoin'
Set-TextCell $ws "C26" 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextCell $ws "D26" '1.514'
Set-TextCell $ws "E26" '  -0.42%  '

Set-TextCell $ws "B27" 'EthereumClassic'
Set-TextCell $ws "C27" 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextCell $ws "D27" '15.30'
Set-TextCell $ws "E27" '  +0.90%  '

Set-TextCell $ws "B28" 'LidoDAOToken'
Set-TextCell $ws "C28" 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextCell $ws "D28" '1.764'
Set-TextCell $ws "E28" '  -1.55%  '

Set-TextCell $ws "B29" 'BitcoinCash'
Set-TextCell $ws "C29" 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextCell $ws "D29" '105.41'
Set-TextCell $ws "E29" '  +2.50%  '

Set-TextCell $ws "B30" 'InternetComputer(DFINITY)'
Set-TextCell $ws "C30" 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextCell $ws "D30" '3.953'
Set-TextCell $ws "E30" '  +5.35%  '

Set-TextCell $ws "B31" 'Stellar'
Set-TextCell $ws "C31" 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextCell $ws "D31" '0.08256'
Set-TextCell $ws "E31" '  -0.69%  '

Set-TextCell $ws "B32" 'Filecoin'
Set-TextCell $ws "C32" 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextCell $ws "D32" '3.653'
Set-TextCell $ws "E32" '  +3.81%  '

Set-TextCell $ws "B33" 'Hedera'
Set-TextCell $ws "C33" 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextCell $ws "D33" '0.04663'
Set-TextCell $ws "E33" '  +3.31%  '

Set-TextCell $ws "B34" 'HuobiToken'
Set-TextCell $ws "C34" 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextCell $ws "D34" '2.652'
Set-TextCell $ws "E34" '  +1.53%  '

Set-TextCell $ws "B35" 'ARBITRUM'
Set-TextCell $ws "C35" 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextCell $ws "D35" '0.9939'
Set-TextCell $ws "E35" '  +1.19%  '

Set-TextCell $ws "B36" 'ImmutableX'
Set-TextCell $ws "C36" 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextCell $ws "D36" '0.6194'
Set-TextCell $ws "E36" '  -0.97%  '

Set-TextCell $ws "B37" 'MXToken'
Set-TextCell $ws "C37" 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextCell $ws "D37" '2.705'
Set-TextCell $ws "E37" '  +0.61%  '

Set-TextCell $ws "B38" 'VeChain'
Set-TextCell $ws "C38" 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextCell $ws "D38" '0.01594'
Set-TextCell $ws "E38" '  +0.20%  '

Set-TextCell $ws "B39" 'RenderToken'
Set-TextCell $ws "C39" 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextCell $ws "D39" '1.915'
Set-TextCell $ws "E39" '  -1.57%  '

Set-TextCell $ws "B40" 'PaxDollar'
Set-TextCell $ws "C40" 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextCell $ws "D40" '0.9989'
Set-TextCell $ws "E40" '  -0.04%  '

Set-TextCell $ws "B41" 'Quant'
Set-TextCell $ws "C41" 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextCell $ws "D41" '100.61'
Set-TextCell $ws "E41" '  +0.72%  '

Set-TextCell $ws "B42" 'TrustWalletToken'
Set-TextCell $ws "C42" 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextCell $ws "D42" '0.7564'
Set-TextCell $ws "E42" '  +2.87%  '

Set-TextCell $ws "B43" 'TheSandbox'
Set-TextCell $ws "C43" 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextCell $ws "D43" '0.3836'
Set-TextCell $ws "E43" '  -1.10%  '

Set-TextCell $ws "B44" 'FraxShare'
Set-TextCell $ws "C44" 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextCell $ws "D44" '5.006'
Set-TextCell $ws "E44" '  -0.58%  '

Set-TextCell $ws "B45" 'Algorand'
Set-TextCell $ws "C45" 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextCell $ws "D45" '0.1129'
Set-TextCell $ws "E45" '  +0.28%  '

Set-TextCell $ws "B46" 'Aptos'
Set-TextCell $ws "C46" 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextCell $ws "D46" '6.261'
Set-TextCell $ws "E46" '  -0.44%  '

Set-TextCell $ws "B47" 'Aave'
Set-TextCell $ws "C47" 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextCell $ws "D47" '55.10'
Set-TextCell $ws "E47" '  +2.67%  '

Set-TextCell $ws "B48" 'Cronos'
Set-TextCell $ws "C48" 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextCell $ws "D48" '0.05219'
Set-TextCell $ws "E48" '  -2.22%  '

Set-TextCell $ws "B49" 'Elrond'
Set-TextCell $ws "C49" 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
Set-TextCell $ws "D49" '30.52'
Set-TextCell $ws "E49" '  +0.93%  '

Set-TextCell $ws "B50" 'EnergySwap'
Set-TextCell $ws "C50" 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell $ws "D50" '7.601'
Set-TextCell $ws "E50" '  -0.12%  '

Set-TextCell $ws "B51" 'Decentraland'
Set-TextCell $ws "C51" 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextCell $ws "D51" '0.3403'
Set-TextCell $ws "E51" '  -1.10%  '
